# Swap the presentation's "Integral" theme palette for the stock
# "Office Theme" palette -- mirroring the author's edit, which
# physically exchanged the contents of ppt/theme/theme1.xml (the theme
# bound to the slide master, i.e. the deck's real design) and
# ppt/theme/theme2.xml (the theme bound only to the notes master).
#
# This host's object model resolves every Theme/ThemeColorScheme
# accessor (SlideMaster.Theme, NotesMaster.Theme, HandoutMaster.Theme,
# Slide.ThemeColorScheme, ...) to the single design theme part
# (ppt/theme/theme1.xml), so that's the part these writes land on.
#
# RGB() isn't available in this host, so colours are supplied as the
# plain decimal BGR integer PowerPoint stores on ColorFormat.RGB
# (0x00BBGGRR), i.e. R | (G << 8) | (B << 16).

$p = $ppt.ActivePresentation

$theme = $p.SlideMaster.Theme
$cs = $theme.ThemeColorScheme

$cs.Colors(1).RGB  = 0          # dk1      000000
$cs.Colors(2).RGB  = 16777215   # lt1      FFFFFF
$cs.Colors(3).RGB  = 6968388    # dk2      44546A
$cs.Colors(4).RGB  = 15132391   # lt2      E7E6E6
$cs.Colors(5).RGB  = 13998939   # accent1  5B9BD5
$cs.Colors(6).RGB  = 3243501    # accent2  ED7D31
$cs.Colors(7).RGB  = 10855845   # accent3  A5A5A5
$cs.Colors(8).RGB  = 49407      # accent4  FFC000
$cs.Colors(9).RGB  = 12874308   # accent5  4472C4
$cs.Colors(10).RGB = 4697456    # accent6  70AD47
$cs.Colors(11).RGB = 12673797   # hlink    0563C1
$cs.Colors(12).RGB = 7491477    # folHlink 954F72

# Best-effort: real PowerPoint also renames the theme/colour-scheme
# when a stock theme is applied. This host does not persist these (no
# file-system / template import), but setting them is harmless if
# unsupported.
try { $theme.Name = "Office Theme" } catch {}
try { $cs.Name = "Office" } catch {}
